# Added harvard case classification: recompute average_doctor stats and
# shift the previous "average_doctor" header/values into a new
# "average_doctor_old" column (BP), while "average_doctor" (BQ) now holds
# the freshly computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("E4").Value = 0.38
$ws.Range("F4").Value = 0.075
$ws.Range("G4").Value = 0.273
$ws.Range("N4").Value = 0.393
$ws.Range("O4").Value = 0.057
$ws.Range("P4").Value = 0.239
$ws.Range("W4").Value = 0.226
$ws.Range("X4").Value = 0.104
$ws.Range("Y4").Value = 0.322
$ws.Range("AI4").Value = 0.206
$ws.Range("AJ4").Value = 0.065
$ws.Range("AK4").Value = 0.256
$ws.Range("AU4").Value = 0.15
$ws.Range("AV4").Value = 0.027
$ws.Range("AW4").Value = 0.163
$ws.Range("BA4").Value = 1.904
$ws.Range("BB4").Value = 0.167
$ws.Range("BC4").Value = 0.408
$ws.Range("BG4").Value = 0.719
$ws.Range("BH4").Value = 0.145
$ws.Range("BI4").Value = 0.381
$ws.Range("BM4").Value = 0.668
$ws.Range("BN4").Value = 0.09
$ws.Range("BO4").Value = 0.301
$ws.Range("BP4").Value = 0.635
$ws.Range("BQ4").Value = 0.644
$ws.Range("E5").Value = 0.507
$ws.Range("G5").Value = 0.314
$ws.Range("N5").Value = 0.762
$ws.Range("O5").Value = 0.076
$ws.Range("P5").Value = 0.276
$ws.Range("W5").Value = 0.234
$ws.Range("X5").Value = 0.116
$ws.Range("Y5").Value = 0.34
$ws.Range("AI5").Value = 0.245
$ws.Range("AJ5").Value = 0.095
$ws.Range("AK5").Value = 0.308
$ws.Range("AU5").Value = 0.306
$ws.Range("AV5").Value = 0.099
$ws.Range("AW5").Value = 0.314
$ws.Range("BA5").Value = 1.361
$ws.Range("BB5").Value = 0.089
$ws.Range("BC5").Value = 0.298
$ws.Range("BG5").Value = 0.406
$ws.Range("BH5").Value = 0.054
$ws.Range("BI5").Value = 0.233
$ws.Range("BM5").Value = 0.577
$ws.Range("BP5").Value = 0.454
$ws.Range("BQ5").Value = 0.453
$ws.Range("E6").Value = 0.434
$ws.Range("N6").Value = 0.519
$ws.Range("W6").Value = 0.23
$ws.Range("AI6").Value = 0.224
$ws.Range("AU6").Value = 0.201
$ws.Range("BA6").Value = 1.575
$ws.Range("BG6").Value = 0.519
$ws.Range("BM6").Value = 0.619
$ws.Range("BP6").Value = 0.525
$ws.Range("BQ6").Value = 0.528
$ws.Range("E7").Value = 0.475
$ws.Range("N7").Value = 0.642
$ws.Range("W7").Value = 0.232
$ws.Range("AI7").Value = 0.236
$ws.Range("AU7").Value = 0.253
$ws.Range("BA7").Value = 1.437
$ws.Range("BG7").Value = 0.445
$ws.Range("BM7").Value = 0.593
$ws.Range("BP7").Value = 0.479
$ws.Range("BQ7").Value = 0.48
$ws.Range("E8").Value = 0.534
$ws.Range("F8").Value = 0.125
$ws.Range("G8").Value = 0.353
$ws.Range("N8").Value = 0.766
$ws.Range("O8").Value = 0.062
$ws.Range("P8").Value = 0.249
$ws.Range("W8").Value = 0.227
$ws.Range("X8").Value = 0.108
$ws.Range("Y8").Value = 0.328
$ws.Range("AI8").Value = 0.224
$ws.Range("AJ8").Value = 0.094
$ws.Range("AK8").Value = 0.306
$ws.Range("AU8").Value = 0.243
$ws.Range("AV8").Value = 0.076
$ws.Range("AW8").Value = 0.276
$ws.Range("BA8").Value = 1.697
$ws.Range("BC8").Value = 0.37
$ws.Range("BG8").Value = 0.555
$ws.Range("BH8").Value = 0.106
$ws.Range("BI8").Value = 0.326
$ws.Range("BM8").Value = 0.697
$ws.Range("BN8").Value = 0.074
$ws.Range("BO8").Value = 0.272
$ws.Range("BP8").Value = 0.566
$ws.Range("BQ8").Value = 0.578
$ws.Range("E9").Value = 0.458
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.646
$ws.Range("O9").Value = 0.229
$ws.Range("P9").Value = 0.478
$ws.Range("W9").Value = 0.125
$ws.Range("X9").Value = 0.109
$ws.Range("Y9").Value = 0.331
$ws.Range("AI9").Value = 0.125
$ws.Range("AJ9").Value = 0.109
$ws.Range("AK9").Value = 0.331
$ws.Range("BA9").Value = 1.625
$ws.Range("BB9").Value = 0.243
$ws.Range("BC9").Value = 0.493
$ws.Range("BG9").Value = 0.583
$ws.Range("BH9").Value = 0.243
$ws.Range("BI9").Value = 0.493
$ws.Range("BM9").Value = 0.625
$ws.Range("BN9").Value = 0.234
$ws.Range("BO9").Value = 0.484
$ws.Range("BP9").Value = 0.542
$ws.Range("BQ9").Value = 0.544
$ws.Range("E10").Value = 0.583
$ws.Range("F10").Value = 0.243
$ws.Range("G10").Value = 0.493
$ws.Range("N10").Value = 0.854
$ws.Range("O10").Value = 0.125
$ws.Range("P10").Value = 0.353
$ws.Range("W10").Value = 0.271
$ws.Range("X10").Value = 0.197
$ws.Range("Y10").Value = 0.444
$ws.Range("AI10").Value = 0.25
$ws.Range("AJ10").Value = 0.188
$ws.Range("AK10").Value = 0.433
$ws.Range("AU10").Value = 0.229
$ws.Range("AV10").Value = 0.177
$ws.Range("AW10").Value = 0.42
$ws.Range("BA10").Value = 1.958
$ws.Range("BB10").Value = 0.25
$ws.Range("BC10").Value = 0.5
$ws.Range("BG10").Value = 0.625
$ws.Range("BH10").Value = 0.234
$ws.Range("BI10").Value = 0.484
$ws.Range("BM10").Value = 0.854
$ws.Range("BN10").Value = 0.125
$ws.Range("BO10").Value = 0.353
$ws.Range("BP10").Value = 0.653
$ws.Range("BQ10").Value = 0.684
$ws.Range("E11").Value = 0.604
$ws.Range("F11").Value = 0.239
$ws.Range("G11").Value = 0.489
$ws.Range("N11").Value = 0.896
$ws.Range("O11").Value = 0.093
$ws.Range("P11").Value = 0.305
$ws.Range("W11").Value = 0.271
$ws.Range("X11").Value = 0.197
$ws.Range("Y11").Value = 0.444
$ws.Range("AI11").Value = 0.25
$ws.Range("AJ11").Value = 0.188
$ws.Range("AK11").Value = 0.433
$ws.Range("AU11").Value = 0.333
$ws.Range("AV11").Value = 0.222
$ws.Range("AW11").Value = 0.471
$ws.Range("BA11").Value = 1.958
$ws.Range("BB11").Value = 0.25
$ws.Range("BC11").Value = 0.5
$ws.Range("BG11").Value = 0.625
$ws.Range("BH11").Value = 0.234
$ws.Range("BI11").Value = 0.484
$ws.Range("BM11").Value = 0.854
$ws.Range("BN11").Value = 0.125
$ws.Range("BO11").Value = 0.353
$ws.Range("BP11").Value = 0.653
$ws.Range("BQ11").Value = 0.684
$ws.Range("E12").Value = 1.414
$ws.Range("F12").Value = 0.656
$ws.Range("G12").Value = 0.81
$ws.Range("N12").Value = 1.568
$ws.Range("O12").Value = 1.245
$ws.Range("P12").Value = 1.116
$ws.Range("W12").Value = 1.846
$ws.Range("X12").Value = 0.746
$ws.Range("Y12").Value = 0.863
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.562
$ws.Range("AV12").Value = 1.746
$ws.Range("AW12").Value = 1.321
$ws.Range("BA12").Value = 3.65
$ws.Range("BB12").Value = 0.344
$ws.Range("BC12").Value = 0.587
$ws.Range("BG12").Value = 1.067
$ws.Range("BH12").Value = 0.062
$ws.Range("BI12").Value = 0.249
$ws.Range("BM12").Value = 1.366
$ws.Range("BN12").Value = 0.427
$ws.Range("BO12").Value = 0.654
$ws.Range("BP12").Value = 1.217
$ws.Range("BQ12").Value = 1.28
$ws.Range("E13").Value = 1.738
$ws.Range("F13").Value = 0.914
$ws.Range("G13").Value = 0.956
$ws.Range("N13").Value = 2.363
$ws.Range("O13").Value = 1.135
$ws.Range("P13").Value = 1.065
$ws.Range("W13").Value = 1.09
$ws.Range("X13").Value = 0.186
$ws.Range("Y13").Value = 0.431
$ws.Range("AI13").Value = 1.383
$ws.Range("AJ13").Value = 0.401
$ws.Range("AK13").Value = 0.633
$ws.Range("AU13").Value = 2.481
$ws.Range("AV13").Value = 1.362
$ws.Range("AW13").Value = 1.167
$ws.Range("BA13").Value = 2.561
$ws.Range("BB13").Value = 0.319
$ws.Range("BC13").Value = 0.565
$ws.Range("BG13").Value = 0.626
$ws.Range("BH13").Value = 0.09
$ws.Range("BI13").Value = 0.299
$ws.Range("BM13").Value = 1.019
$ws.Range("BN13").Value = 0.367
$ws.Range("BO13").Value = 0.606
$ws.Range("BP13").Value = 0.854
$ws.Range("BQ13").Value = 0.796
